$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy header style from C1 into new header cells D1:E1
$ws.Range("C1").Copy() | Out-Null
$ws.Range("D1:E1").PasteSpecial(-4122) | Out-Null

# Header row
$ws.Range("D1").Value = 3
$ws.Range("E1").Value = 4

$ws.Range("B2").Value = -0.3171415177393752
$ws.Range("C2").Value = -0.3040417462236097
$ws.Range("D2").Value = -0.291644377748045
$ws.Range("E2").Value = -0.2803157617148355
$ws.Range("B3").Value = 0.2251170848256619
$ws.Range("C3").Value = 0.2251535355371952
$ws.Range("D3").Value = 0.2258924693243424
$ws.Range("E3").Value = 0.2268734438303079
$ws.Range("B4").Value = 0.1456377525183181
$ws.Range("C4").Value = 0.1544151644258079
$ws.Range("D4").Value = 0.1634582775498603
$ws.Range("E4").Value = 0.1721354018535958
$ws.Range("B5").Value = -0.07508467623793266
$ws.Range("C5").Value = -0.06922912467150039
$ws.Range("D5").Value = -0.06389751015466108
$ws.Range("E5").Value = -0.05931101797720713
$ws.Range("B6").Value = 0.1508541478869574
$ws.Range("C6").Value = 0.1592649206024101
$ws.Range("D6").Value = 0.1666633624084359
$ws.Range("E6").Value = 0.1728513362613026
$ws.Range("B7").Value = -0.455675197302285
$ws.Range("C7").Value = -0.4466575725286778
$ws.Range("D7").Value = -0.4374948230304833
$ws.Range("E7").Value = -0.428624097317399
$ws.Range("B8").Value = -0.2675861058738578
$ws.Range("C8").Value = -0.2530181817640409
$ws.Range("D8").Value = -0.2387230257345911
$ws.Range("E8").Value = -0.2252340431466514
$ws.Range("B9").Value = -0.3908074223347818
$ws.Range("C9").Value = -0.3885918744649107
$ws.Range("D9").Value = -0.3847800886814022
$ws.Range("E9").Value = -0.3799580678948297
$ws.Range("B10").Value = 0.3596325379356408
$ws.Range("C10").Value = 0.3690991585549788
$ws.Range("D10").Value = 0.377416046792225
$ws.Range("E10").Value = 0.3843099641951215
$ws.Range("B11").Value = -0.2466499125343932
$ws.Range("C11").Value = -0.2374145889386449
$ws.Range("D11").Value = -0.2283606933100192
$ws.Range("E11").Value = -0.2199436261761657
$ws.Range("B12").Value = -0.07800631489799899
$ws.Range("C12").Value = -0.06832904299805047
$ws.Range("D12").Value = -0.06072853088561912
$ws.Range("E12").Value = -0.05501965891714022
$ws.Range("B13").Value = -0.008176485672542824
$ws.Range("C13").Value = -0.01024573843619881
$ws.Range("D13").Value = -0.01171705321852328
$ws.Range("E13").Value = -0.01295245704035458
$ws.Range("B14").Value = 0.08665506124432008
$ws.Range("C14").Value = 0.09767730266494998
$ws.Range("D14").Value = 0.105062170142916
$ws.Range("E14").Value = 0.1095508428675377
$ws.Range("B15").Value = 0.02695589753023021
$ws.Range("C15").Value = 0.03554015468291617
$ws.Range("D15").Value = 0.03950052157395213
$ws.Range("E15").Value = 0.04001126288311345
$ws.Range("B16").Value = 0.4036244154493354
$ws.Range("C16").Value = 0.4182123011612361
$ws.Range("D16").Value = 0.426667308028124
$ws.Range("E16").Value = 0.4304674828948816
$ws.Range("B17").Value = 0.6488478964810216
$ws.Range("C17").Value = 0.6463634593515115
$ws.Range("D17").Value = 0.6397009750300002
$ws.Range("E17").Value = 0.6301036481976009
$ws.Range("B18").Value = 0.07863562390646171
$ws.Range("C18").Value = 0.050867864187726
$ws.Range("D18").Value = 0.02598052180645262
$ws.Range("E18").Value = 0.003800673892381534
$ws.Range("B19").Value = 0.4195822116113737
$ws.Range("C19").Value = 0.4182019070579845
$ws.Range("D19").Value = 0.4154713097687464
$ws.Range("E19").Value = 0.4115980141863978
$ws.Range("B20").Value = 0.3467898443532377
$ws.Range("C20").Value = 0.3586001354007108
$ws.Range("D20").Value = 0.3642587293248287
$ws.Range("E20").Value = 0.365298131682758
$ws.Range("B21").Value = 0.5337064219450647
$ws.Range("C21").Value = 0.5605837047146707
$ws.Range("D21").Value = 0.5795560394245626
$ws.Range("E21").Value = 0.5920223842420665
$ws.Range("B22").Value = 0.4192817319970564
$ws.Range("C22").Value = 0.4216504564096156
$ws.Range("D22").Value = 0.4198413741633195
$ws.Range("E22").Value = 0.4149427765125479
$ws.Range("B23").Value = -0.03618381376271146
$ws.Range("C23").Value = -0.02984816026661184
$ws.Range("D23").Value = -0.02762348792757056
$ws.Range("E23").Value = -0.02851268401638681
$ws.Range("B24").Value = 4.326030357861844
$ws.Range("C24").Value = 4.319056669188416
$ws.Range("D24").Value = 4.253385329351993
$ws.Range("E24").Value = 4.144565729608261
$ws.Range("B25").Value = 0.492050565015991
$ws.Range("C25").Value = 0.4598891472662188
$ws.Range("D25").Value = 0.4331238609730279
$ws.Range("E25").Value = 0.4105244079658265
$ws.Range("B26").Value = 0.3860257590419049
$ws.Range("C26").Value = 0.3670980678865449
$ws.Range("D26").Value = 0.3487579803555776
$ws.Range("E26").Value = 0.331080465229419
$ws.Range("B27").Value = 0.2942499834532983
$ws.Range("C27").Value = 0.2621466084748486
$ws.Range("D27").Value = 0.234674336632878
$ws.Range("E27").Value = 0.2102714291201903
$ws.Range("B28").Value = 1.030639423042036
$ws.Range("C28").Value = 0.991487030393535
$ws.Range("D28").Value = 0.9610564509519863
$ws.Range("E28").Value = 0.9354775223652831
$ws.Range("B29").Value = 5.686085270340403
$ws.Range("C29").Value = 5.224988178542493
$ws.Range("D29").Value = 4.79906475202811
$ws.Range("E29").Value = 4.409204672351636
$ws.Range("B30").Value = 0.965142905670423
$ws.Range("C30").Value = 0.923784986971142
$ws.Range("D30").Value = 0.8877063190965784
$ws.Range("E30").Value = 0.855988279479129
$ws.Range("B31").Value = -0.2246992233461589
$ws.Range("C31").Value = -0.2707495582110871
$ws.Range("D31").Value = -0.3067469171251993
$ws.Range("E31").Value = -0.3354843502624525
$ws.Range("B32").Value = 0.7843652602154243
$ws.Range("C32").Value = 0.7504266855551888
$ws.Range("D32").Value = 0.7220160779612326
$ws.Range("E32").Value = 0.6974390202878871
$ws.Range("B33").Value = 0.9100058871148673
$ws.Range("C33").Value = 0.8844208549016394
$ws.Range("D33").Value = 0.863580810902531
$ws.Range("E33").Value = 0.8459464034106425
$ws.Range("B34").Value = -0.6797048748884066
$ws.Range("C34").Value = -0.7032468508123622
$ws.Range("D34").Value = -0.7218037139361796
$ws.Range("E34").Value = -0.7368140142040426
$ws.Range("B35").Value = 0.8008466521657421
$ws.Range("C35").Value = 0.800343452240891
$ws.Range("D35").Value = 0.8008870183887973
$ws.Range("E35").Value = 0.8015653773189956
$ws.Range("B36").Value = 0.7522427518756689
$ws.Range("C36").Value = 0.7483092607616741
$ws.Range("D36").Value = 0.7467920166938515
$ws.Range("E36").Value = 0.7464495457261444
$ws.Range("B37").Value = 0.725255320649142
$ws.Range("C37").Value = 0.7198224738286401
$ws.Range("D37").Value = 0.7170416746796945
$ws.Range("E37").Value = 0.71564121926974
$ws.Range("B38").Value = 0.7189995198482038
$ws.Range("C38").Value = 0.7053106754767264
$ws.Range("D38").Value = 0.6944415201437384
$ws.Range("E38").Value = 0.6853932015980665
$ws.Range("B39").Value = 0.5663840323053958
$ws.Range("C39").Value = 0.5701106224081278
$ws.Range("D39").Value = 0.5745964938653605
$ws.Range("E39").Value = 0.5789574463301435
$ws.Range("B40").Value = 0.7343382325555655
$ws.Range("C40").Value = 0.7401332273910494
$ws.Range("D40").Value = 0.7459332016200764
$ws.Range("E40").Value = 0.7509735553210193
$ws.Range("B41").Value = 0.5550826488079144
$ws.Range("C41").Value = 0.5508581871494762
$ws.Range("D41").Value = 0.5486856884849532
$ws.Range("E41").Value = 0.5475857123802477
$ws.Range("B42").Value = 0.6771058842947451
$ws.Range("C42").Value = 0.6645571692077971
$ws.Range("D42").Value = 0.6554664128341929
$ws.Range("E42").Value = 0.6485338827807091
$ws.Range("B43").Value = 0.7073435764176282
$ws.Range("C43").Value = 0.6986510744487012
$ws.Range("D43").Value = 0.6923441292037344
$ws.Range("E43").Value = 0.6874174817595955
$ws.Range("B44").Value = 0.6637697425569724
$ws.Range("C44").Value = 0.6638463159649135
$ws.Range("D44").Value = 0.6655216730290923
$ws.Range("E44").Value = 0.6677474596245438
$ws.Range("B45").Value = 0.6303503093024357
$ws.Range("C45").Value = 0.6316393742209037
$ws.Range("D45").Value = 0.6355093236849141
$ws.Range("E45").Value = 0.6404919044044228
$ws.Range("B46").Value = -1.28606498749736
$ws.Range("C46").Value = -1.281318576789337
$ws.Range("D46").Value = -1.27586495541497
$ws.Range("E46").Value = -1.270040349946611
$ws.Range("B47").Value = -0.9985240308883783
$ws.Range("C47").Value = -0.9953210513289173
$ws.Range("D47").Value = -0.9912822502027138
$ws.Range("E47").Value = -0.9868107119277745
$ws.Range("B48").Value = -0.8931648659201764
$ws.Range("C48").Value = -0.8891017246582267
$ws.Range("D48").Value = -0.8836044856295788
$ws.Range("E48").Value = -0.8773203593483833
$ws.Range("B49").Value = -0.6587458580148532
$ws.Range("C49").Value = -0.6541934193054134
$ws.Range("D49").Value = -0.6486546449944572
$ws.Range("E49").Value = -0.6427025874187626
$ws.Range("B50").Value = -0.05851859617447983
$ws.Range("C50").Value = -0.05623713313460828
$ws.Range("D50").Value = -0.05347623463888525
$ws.Range("E50").Value = -0.05075959431968712
$ws.Range("B51").Value = -0.8878024941345092
$ws.Range("C51").Value = -0.88181270688059
$ws.Range("D51").Value = -0.8747617565870891
$ws.Range("E51").Value = -0.8672481028329181
$ws.Range("B52").Value = -0.8878024941345092
$ws.Range("C52").Value = -0.88181270688059
$ws.Range("D52").Value = -0.8747617565870891
$ws.Range("E52").Value = -0.8672481028329181
$ws.Range("B53").Value = -1.114620358715956
$ws.Range("C53").Value = -1.116975934909187
$ws.Range("D53").Value = -1.116815806966988
$ws.Range("E53").Value = -1.114871712192338
$ws.Range("B54").Value = -0.2008549705718582
$ws.Range("C54").Value = -0.1947372290805855
$ws.Range("D54").Value = -0.1882734008913984
$ws.Range("E54").Value = -0.1820128069211368
$ws.Range("B55").Value = -1.0192210989668
$ws.Range("C55").Value = -1.015048157973971
$ws.Range("D55").Value = -1.010221175374169
$ws.Range("E55").Value = -1.00510704194382
$ws.Range("B56").Value = -0.9080489991024606
$ws.Range("C56").Value = -0.8977876783020756
$ws.Range("D56").Value = -0.8885832658536043
$ws.Range("E56").Value = -0.8804858252887116
$ws.Range("B57").Value = -0.9381462945588478
$ws.Range("C57").Value = -0.9300980938665987
$ws.Range("D57").Value = -0.9232693827135859
$ws.Range("E57").Value = -0.9174869799337892
$ws.Range("B58").Value = -1.142081707100022
$ws.Range("C58").Value = -1.123637583142322
$ws.Range("D58").Value = -1.107984300774948
$ws.Range("E58").Value = -1.094755940194904
$ws.Range("B59").Value = -0.866513041817719
$ws.Range("C59").Value = -0.8511685513275964
$ws.Range("D59").Value = -0.836862266073065
$ws.Range("E59").Value = -0.8237826303203279
$ws.Range("B60").Value = -0.5075653945589996
$ws.Range("C60").Value = -0.4902850983854158
$ws.Range("D60").Value = -0.4751437649096459
$ws.Range("E60").Value = -0.4621057241465001
$ws.Range("B61").Value = 0.3692306399011779
$ws.Range("C61").Value = 0.3716943380942198
$ws.Range("D61").Value = 0.3745374886554746
$ws.Range("E61").Value = 0.377174682581613
$ws.Range("B62").Value = -1.217607622281926
$ws.Range("C62").Value = -1.203251063982498
$ws.Range("D62").Value = -1.191225051707048
$ws.Range("E62").Value = -1.181110987235549
$ws.Range("B63").Value = -0.7624594456434794
$ws.Range("C63").Value = -0.7335006093035624
$ws.Range("D63").Value = -0.7068733547044074
$ws.Range("E63").Value = -0.6828154500169004
$ws.Range("B64").Value = -0.9044845088990647
$ws.Range("C64").Value = -0.8991233070265447
$ws.Range("D64").Value = -0.892690924503508
$ws.Range("E64").Value = -0.8856952084076489
$ws.Range("B65").Value = -0.1236732446001867
$ws.Range("C65").Value = -0.1057560117175542
$ws.Range("D65").Value = -0.0900461125683942
$ws.Range("E65").Value = -0.0766328719526629
$ws.Range("B66").Value = -0.8025660841616801
$ws.Range("C66").Value = -0.7847113697185064
$ws.Range("D66").Value = -0.7703223723821625
$ws.Range("E66").Value = -0.7588869234675441
$ws.Range("B67").Value = -0.7845681523996909
$ws.Range("C67").Value = -0.7583261251229748
$ws.Range("D67").Value = -0.7377665429980237
$ws.Range("E67").Value = -0.7219602109232565
